# Apply the changes described by the commit "Push Anh Hoan Push"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) NganHang sheet: update a few balance figures
# ---------------------------------------------------------------------
$wsBank = $wb.Worksheets.Item("NganHang")
$wsBank.Range("B3").Value = 360420000
$wsBank.Range("B5").Value = 43870000
$wsBank.Range("B12").Value = 7970000

# ---------------------------------------------------------------------
# 2) XeMay sheet: add a new motorbike entry in row 13
# ---------------------------------------------------------------------
$wsMoto = $wb.Worksheets.Item("XeMay")
$wsMoto.Range("A13").Value = 123
$wsMoto.Range("B13").Value = "Dream2018"
$wsMoto.Range("D13").Value = 100
$wsMoto.Range("E13").Value = "Không"
$wsMoto.Range("F13").Value = "Du lịch"
$wsMoto.Range("G13").Value = 1000000
$wsMoto.Range("H13").Value = 100000
$wsMoto.Range("I13").Value = 100000
$wsMoto.Range("J13").Value = 100000
$wsMoto.Range("K13").Value = 100000
$wsMoto.Range("L13").Value = 100000
$wsMoto.Range("M13").Value = 100000
$wsMoto.Range("N13").Value = "47-U1"
$wsMoto.Range("O13").Value = "Sai"

# Keep the same date formatting as the rows above it, then set the date
$wsMoto.Range("C12").Copy()
$wsMoto.Range("C13").PasteSpecial(-4122)
$wsMoto.Range("C13").Value = 45265

# Column B got a bit wider to fit the new model name
$wsMoto.Columns.Item(2).ColumnWidth = 9.6

# Update the view: scroll back to the left and move the selection
$wsMoto.Activate()
$wsMoto.Range("F19").Select()
